# This workbook's data rows got re-ordered upstream (records 12/13, 19-23,
# and 31-35 were re-sequenced) while keeping the same physical spreadsheet
# row numbers. Re-create that by rotating the row *contents* between the
# fixed row positions.
#
# Columns X, Y, Z, AA, AB (Externid/Startdatum/Starttid/Slutdatum/Sluttid)
# are identical across every row in these blocks, so they are left alone -
# we only move A:W and AB:AY ... actually AB:AY covers AC.. through AY, and
# A:W covers A through W. Y/Z/AA/AB are skipped entirely since they never
# differ between any of the source/destination rows involved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSnapshot {
    param($ws, $row)
    $snap = @{}
    $snap["Part1"] = $ws.Range("A$row`:W$row").Value()
    $snap["Part2"] = $ws.Range("AB$row`:AY$row").Value()
    return $snap
}

function Set-RowSnapshot {
    param($ws, $row, $snap)
    $ws.Range("A$row`:W$row").Value = $snap["Part1"]
    $ws.Range("AB$row`:AY$row").Value = $snap["Part2"]
}

# --- Block 1: rows 12 <-> 13 (simple swap) ---
$r12 = Get-RowSnapshot $ws 12
$r13 = Get-RowSnapshot $ws 13
Set-RowSnapshot $ws 12 $r13
Set-RowSnapshot $ws 13 $r12

# --- Block 2: rows 19,20,21,22,23 cyclic rotation ---
# new(19)=old(23); new(20)=old(21); new(21)=old(22); new(22)=old(19); new(23)=old(20)
$r19 = Get-RowSnapshot $ws 19
$r20 = Get-RowSnapshot $ws 20
$r21 = Get-RowSnapshot $ws 21
$r22 = Get-RowSnapshot $ws 22
$r23 = Get-RowSnapshot $ws 23

Set-RowSnapshot $ws 19 $r23
Set-RowSnapshot $ws 20 $r21
Set-RowSnapshot $ws 21 $r22
Set-RowSnapshot $ws 22 $r19
Set-RowSnapshot $ws 23 $r20

# --- Block 3: rows 31,32,33,34,35 cyclic rotation ---
# new(31)=old(32); new(32)=old(33); new(33)=old(34); new(34)=old(35); new(35)=old(31)
$r31 = Get-RowSnapshot $ws 31
$r32 = Get-RowSnapshot $ws 32
$r33 = Get-RowSnapshot $ws 33
$r34 = Get-RowSnapshot $ws 34
$r35 = Get-RowSnapshot $ws 35

Set-RowSnapshot $ws 31 $r32
Set-RowSnapshot $ws 32 $r33
Set-RowSnapshot $ws 33 $r34
Set-RowSnapshot $ws 34 $r35
Set-RowSnapshot $ws 35 $r31

Write-Host "Row rotations applied."
